# Auto-generated edit script applying cached market-data value updates
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(41, 8).Value = 946.65
$ws.Cells.Item(41, 9).Value = 1560
$ws.Cells.Item(41, 10).Value = 683.7857
$ws.Cells.Item(41, 11).Value = 1560
$ws.Cells.Item(41, 12).Value = 683.7857
$ws.Cells.Item(41, 13).Value = -1120
$ws.Cells.Item(41, 14).Value = -1563.7857
$ws.Cells.Item(76, 8).Value = 4687085
$ws.Cells.Item(79, 8).Value = 4687085
$ws.Cells.Item(80, 8).Value = 1408.8
$ws.Cells.Item(80, 9).Value = 1491.5333
$ws.Cells.Item(80, 10).Value = 1160.6
$ws.Cells.Item(80, 11).Value = 4474.5999
$ws.Cells.Item(80, 12).Value = 3481.8
$ws.Cells.Item(80, 13).Value = -3476.5999
$ws.Cells.Item(80, 14).Value = -5477.799999999999
$ws.Cells.Item(83, 8).Value = 1408.8
$ws.Cells.Item(83, 9).Value = 1491.5333
$ws.Cells.Item(83, 10).Value = 1160.6
$ws.Cells.Item(83, 11).Value = 13423.7997
$ws.Cells.Item(83, 12).Value = 10445.4
$ws.Cells.Item(83, 13).Value = -8431.7997
$ws.Cells.Item(83, 14).Value = -20429.4
$ws.Cells.Item(86, 8).Value = 82781.8
$ws.Cells.Item(86, 9).Value = 154618.38
$ws.Cells.Item(86, 10).Value = 682.8570999999999
$ws.Cells.Item(86, 11).Value = 154618.38
$ws.Cells.Item(86, 12).Value = 682.8570999999999
$ws.Cells.Item(86, 13).Value = -153495.38
$ws.Cells.Item(86, 14).Value = -2928.8571
$ws.Cells.Item(88, 8).Value = 3764.7144
$ws.Cells.Item(88, 10).Value = 4468.6
$ws.Cells.Item(88, 12).Value = 4468.6
$ws.Cells.Item(88, 14).Value = -5280.6
$ws.Cells.Item(89, 8).Value = 82781.8
$ws.Cells.Item(89, 9).Value = 154618.38
$ws.Cells.Item(89, 10).Value = 682.8570999999999
$ws.Cells.Item(89, 11).Value = 773091.9
$ws.Cells.Item(89, 12).Value = 3414.2855
$ws.Cells.Item(89, 13).Value = -767475.9
$ws.Cells.Item(89, 14).Value = -14646.2855
$ws.Cells.Item(91, 8).Value = 3764.7144
$ws.Cells.Item(91, 10).Value = 4468.6
$ws.Cells.Item(91, 12).Value = 4468.6
$ws.Cells.Item(91, 14).Value = -7276.6
$ws.Cells.Item(100, 8).Value = 1190.6923
$ws.Cells.Item(100, 9).Value = 1007.2727
$ws.Cells.Item(100, 11).Value = 1007.2727
$ws.Cells.Item(100, 13).Value = -466.2727
$ws.Cells.Item(117, 8).Value = 48871
$ws.Cells.Item(117, 10).Value = 48871
$ws.Cells.Item(117, 12).Value = 48871
$ws.Cells.Item(117, 14).Value = -58049
$ws.Cells.Item(135, 8).Value = 378.91666
$ws.Cells.Item(135, 10).Value = 125
$ws.Cells.Item(135, 12).Value = 1125
$ws.Cells.Item(135, 14).Value = -6195
$ws.Cells.Item(138, 8).Value = 2105.5933
$ws.Cells.Item(138, 9).Value = 1884.0646
$ws.Cells.Item(138, 10).Value = 2350.8572
$ws.Cells.Item(138, 11).Value = 5652.1938
$ws.Cells.Item(138, 12).Value = 7052.571599999999
$ws.Cells.Item(138, 13).Value = -512.1938
$ws.Cells.Item(138, 14).Value = -17332.5716
$ws.Cells.Item(141, 8).Value = 2792.3333
$ws.Cells.Item(141, 9).Value = 1046.7142
$ws.Cells.Item(141, 10).Value = 6283.5713
$ws.Cells.Item(141, 11).Value = 3140.1426
$ws.Cells.Item(141, 12).Value = 18850.7139
$ws.Cells.Item(141, 13).Value = 2039.8574
$ws.Cells.Item(141, 14).Value = -29210.7139

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2282.4885
$ws.Cells.Item(32, 9).Value = 1748.0476
$ws.Cells.Item(32, 11).Value = 1748.0476
$ws.Cells.Item(32, 13).Value = -1461.0476
$ws.Cells.Item(108, 8).Value = 50000
$ws.Cells.Item(108, 10).Value = 50000
$ws.Cells.Item(108, 12).Value = 50000
$ws.Cells.Item(108, 14).Value = -57680
$ws.Cells.Item(122, 8).Value = 1099.8636
$ws.Cells.Item(122, 9).Value = 1198.0625
$ws.Cells.Item(122, 11).Value = 3594.1875
$ws.Cells.Item(122, 13).Value = -1144.1875
$ws.Cells.Item(123, 8).Value = 62000
$ws.Cells.Item(123, 10).Value = 62000
$ws.Cells.Item(123, 12).Value = 62000
$ws.Cells.Item(123, 14).Value = -71800
$ws.Cells.Item(132, 8).Value = 2357.8333
$ws.Cells.Item(132, 9).Value = 1580.3077
$ws.Cells.Item(132, 11).Value = 4740.9231
$ws.Cells.Item(132, 13).Value = -2210.9231

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 185897.9
$ws.Cells.Item(86, 10).Value = 201987.7
$ws.Cells.Item(86, 12).Value = 201987.7
$ws.Cells.Item(86, 14).Value = -204233.7
$ws.Cells.Item(89, 8).Value = 185897.9
$ws.Cells.Item(89, 10).Value = 201987.7
$ws.Cells.Item(89, 12).Value = 1009938.5
$ws.Cells.Item(89, 14).Value = -1021170.5
$ws.Cells.Item(135, 8).Value = 60000
$ws.Cells.Item(135, 10).Value = 60000
$ws.Cells.Item(135, 12).Value = 60000
$ws.Cells.Item(135, 14).Value = -70140

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 1176984.6
$ws.Cells.Item(58, 9).Value = 1554160.9
$ws.Cells.Item(58, 10).Value = 3547.3333
$ws.Cells.Item(58, 11).Value = 1554160.9
$ws.Cells.Item(58, 12).Value = 3547.3333
$ws.Cells.Item(58, 13).Value = -1553957.9
$ws.Cells.Item(58, 14).Value = -3953.3333
$ws.Cells.Item(109, 8).Value = 0
$ws.Cells.Item(109, 10).Value = 0
$ws.Cells.Item(109, 12).Value = 0
$ws.Cells.Item(109, 14).ClearContents()
$ws.Cells.Item(122, 8).Value = 3942.889
$ws.Cells.Item(122, 9).Value = 994.6
$ws.Cells.Item(122, 11).Value = 2983.8
$ws.Cells.Item(122, 13).Value = -533.8000000000002
$ws.Cells.Item(132, 8).Value = 2171.4243
$ws.Cells.Item(132, 9).Value = 1173.1
$ws.Cells.Item(132, 10).Value = 3707.3076
$ws.Cells.Item(132, 11).Value = 3519.3
$ws.Cells.Item(132, 12).Value = 11121.9228
$ws.Cells.Item(132, 13).Value = -989.2999999999997
$ws.Cells.Item(132, 14).Value = -16181.9228
$ws.Cells.Item(136, 8).Value = 1176984.6
$ws.Cells.Item(136, 9).Value = 1554160.9
$ws.Cells.Item(136, 10).Value = 3547.3333
$ws.Cells.Item(136, 11).Value = 4662482.699999999
$ws.Cells.Item(136, 12).Value = 10641.9999
$ws.Cells.Item(136, 13).Value = -4659932.699999999
$ws.Cells.Item(136, 14).Value = -15741.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 692
$ws.Cells.Item(5, 10).Value = 955
$ws.Cells.Item(5, 12).Value = 2865
$ws.Cells.Item(5, 14).Value = -3089
$ws.Cells.Item(92, 8).Value = 317.54544
$ws.Cells.Item(92, 9).Value = 297.5
$ws.Cells.Item(92, 10).Value = 329
$ws.Cells.Item(92, 11).Value = 892.5
$ws.Cells.Item(92, 12).Value = 987
$ws.Cells.Item(92, 13).Value = 355.5
$ws.Cells.Item(92, 14).Value = -3483
$ws.Cells.Item(122, 8).Value = 812.5789
$ws.Cells.Item(122, 10).Value = 1123.7273
$ws.Cells.Item(122, 12).Value = 10113.5457
$ws.Cells.Item(122, 14).Value = -15013.5457
$ws.Cells.Item(131, 8).Value = 7258772.5
$ws.Cells.Item(131, 10).Value = 13127.062
$ws.Cells.Item(131, 12).Value = 39381.186
$ws.Cells.Item(131, 14).Value = -49461.186
$ws.Cells.Item(135, 8).Value = 692
$ws.Cells.Item(135, 10).Value = 955
$ws.Cells.Item(135, 12).Value = 8595
$ws.Cells.Item(135, 14).Value = -13665

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(108, 8).Value = 50000
$ws.Cells.Item(108, 10).Value = 50000
$ws.Cells.Item(108, 12).Value = 50000
$ws.Cells.Item(108, 14).Value = -57680
$ws.Cells.Item(110, 8).Value = 90000
$ws.Cells.Item(110, 10).Value = 90000
$ws.Cells.Item(110, 12).Value = 90000
$ws.Cells.Item(110, 14).Value = -98180
$ws.Cells.Item(113, 8).Value = 1071.4286
$ws.Cells.Item(113, 9).Value = 800
$ws.Cells.Item(113, 10).Value = 1116.6666
$ws.Cells.Item(113, 11).Value = 800
$ws.Cells.Item(113, 12).Value = 1116.6666
$ws.Cells.Item(113, 13).Value = 1370
$ws.Cells.Item(113, 14).Value = -5456.6666
$ws.Cells.Item(126, 8).Value = 2177585.5
$ws.Cells.Item(126, 9).Value = 4632101
$ws.Cells.Item(126, 11).Value = 13896303
$ws.Cells.Item(126, 13).Value = -13893833

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 4848.143
$ws.Cells.Item(16, 9).Value = 4848.143
$ws.Cells.Item(16, 11).Value = 4848.143
$ws.Cells.Item(16, 13).Value = -4678.143
$ws.Cells.Item(61, 8).Value = 2259.0557
$ws.Cells.Item(61, 9).Value = 2025.9286
$ws.Cells.Item(61, 11).Value = 2025.9286
$ws.Cells.Item(61, 13).Value = -1823.9286
$ws.Cells.Item(68, 8).Value = 1922.5
$ws.Cells.Item(68, 9).Value = 1603.7142
$ws.Cells.Item(68, 10).Value = 2666.3333
$ws.Cells.Item(68, 11).Value = 1603.7142
$ws.Cells.Item(68, 12).Value = 2666.3333
$ws.Cells.Item(68, 13).Value = -854.7141999999999
$ws.Cells.Item(68, 14).Value = -4164.3333
$ws.Cells.Item(71, 8).Value = 1922.5
$ws.Cells.Item(71, 9).Value = 1603.7142
$ws.Cells.Item(71, 10).Value = 2666.3333
$ws.Cells.Item(71, 11).Value = 8018.571
$ws.Cells.Item(71, 12).Value = 13331.6665
$ws.Cells.Item(71, 13).Value = -4274.571
$ws.Cells.Item(71, 14).Value = -20819.6665
$ws.Cells.Item(82, 8).Value = 1650.3334
$ws.Cells.Item(82, 9).Value = 1500.5
$ws.Cells.Item(82, 11).Value = 1500.5
$ws.Cells.Item(82, 13).Value = -1139.5
$ws.Cells.Item(85, 8).Value = 1650.3334
$ws.Cells.Item(85, 9).Value = 1500.5
$ws.Cells.Item(85, 11).Value = 1500.5
$ws.Cells.Item(85, 13).Value = -252.5
$ws.Cells.Item(113, 8).Value = 2259.0557
$ws.Cells.Item(113, 9).Value = 2025.9286
$ws.Cells.Item(113, 11).Value = 2025.9286
$ws.Cells.Item(113, 13).Value = 144.0714

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 3835
$ws.Cells.Item(62, 9).Value = 3251
$ws.Cells.Item(62, 10).Value = 5003
$ws.Cells.Item(62, 11).Value = 3251
$ws.Cells.Item(62, 12).Value = 5003
$ws.Cells.Item(62, 13).Value = -2627
$ws.Cells.Item(62, 14).Value = -6251
$ws.Cells.Item(65, 8).Value = 3835
$ws.Cells.Item(65, 9).Value = 3251
$ws.Cells.Item(65, 10).Value = 5003
$ws.Cells.Item(65, 11).Value = 16255
$ws.Cells.Item(65, 12).Value = 25015
$ws.Cells.Item(65, 13).Value = -13135
$ws.Cells.Item(65, 14).Value = -31255
$ws.Cells.Item(113, 8).Value = 516.3333
$ws.Cells.Item(113, 9).Value = 374.64285
$ws.Cells.Item(113, 10).Value = 1012.25
$ws.Cells.Item(113, 11).Value = 1123.92855
$ws.Cells.Item(113, 12).Value = 3036.75
$ws.Cells.Item(113, 13).Value = 1046.07145
$ws.Cells.Item(113, 14).Value = -7376.75
$ws.Cells.Item(123, 8).Value = 47599.582
$ws.Cells.Item(123, 10).Value = 47599.582
$ws.Cells.Item(123, 12).Value = 47599.582
$ws.Cells.Item(123, 14).Value = -57399.582
